# Commit: "tìm kiém điện thoại theo thương hiệu, thêm 4 sản phẩm vào csdl"
# -> Add a new "Samsung: 4" bullet right after the existing "Iphone: 4"
#    bullet in the "Thêm sản phẩm" / "Điện thoại:" list (same list level,
#    same numbering, as a brand-wise phone count).

$d = $word.ActiveDocument

# Locate the paragraph that reads "Iphone: 4" inside the numbered list
# (ilvl 6 / numId 4) used for the per-brand phone counters.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`v")
    if ($t -eq "Iphone: 4" -and $p.Range.ListFormat.ListLevelNumber -eq 7) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a brand-new paragraph right after it; Word clones the
    # paragraph formatting (ListParagraph style + ilvl 6 / numId 4) from
    # the paragraph it is inserted after.
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()

    # Collapse to the (still empty) new paragraph and inject the two text
    # runs via raw OOXML so that "Samsung: " and "4" remain separate
    # <w:r> runs (matching the source edit) instead of being coalesced
    # into a single run the way plain text insertion would do.
    $nr = $newPara.Range
    $nr.Collapse(1)

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:pPr>' +
             '<w:pStyle w:val="ListParagraph"/>' +
             '<w:numPr><w:ilvl w:val="6"/><w:numId w:val="4"/></w:numPr>' +
           '</w:pPr>' +
           '<w:r><w:t xml:space="preserve">Samsung: </w:t></w:r>' +
           '<w:r><w:t>4</w:t></w:r>' +
           '</w:p>'

    $null = $nr.InsertXML($xml)
}
